$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the upstream cryptos.com price refresh.
# Each text-bearing cell is forced to Text format before the write and
# reset to the default "Normal" style afterward so values such as
# "0.180" / "41.15" survive as literal strings instead of being
# reinterpreted (and rounded) as floating point numbers by Excel.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '63.262.30'
Set-TextValue 'E2' '  +5.48%  '

# Row 3
Set-TextValue 'D3' '3.384.56'
Set-TextValue 'E3' '  +6.08%  '

# Row 4
Set-TextValue 'E4' '  +0.02%  '

# Row 5
Set-TextValue 'D5' '577.47'
Set-TextValue 'E5' '  +7.75%  '

# Row 6
Set-TextValue 'D6' '154.40'
Set-TextValue 'E6' '  +6.48%  '

# Row 7
Set-TextValue 'E7' '  +0.12%  '

# Row 8
Set-TextValue 'D8' '3.387.59'
Set-TextValue 'E8' '  +6.00%  '

# Row 9
Set-TextValue 'E9' '  +0.46%  '

# Row 10
Set-TextValue 'D10' '7.47'
Set-TextValue 'E10' '  +2.06%  '

# Row 11
Set-TextValue 'E11' '  +6.57%  '

# Row 12
Set-TextValue 'E12' '  +1.11%  '

# Row 13
Set-TextValue 'D13' '3.972.36'
Set-TextValue 'E13' '  +6.12%  '

# Row 14
Set-TextValue 'E14' '  +0.39%  '

# Row 15
Set-TextValue 'E15' '  +6.84%  '

# Row 16
Set-TextValue 'D16' '26.97'
Set-TextValue 'E16' '  +4.60%  '

# Row 17
Set-TextValue 'D17' '63.377.80'
Set-TextValue 'E17' '  +5.60%  '

# Row 18
Set-TextValue 'D18' '3.390.85'
Set-TextValue 'E18' '  +6.10%  '

# Row 19
Set-TextValue 'D19' '6.37'
Set-TextValue 'E19' '  +1.68%  '

# Row 20
Set-TextValue 'D20' '13.95'
Set-TextValue 'E20' '  +5.29%  '

# Row 21
Set-TextValue 'D21' '8.41'
Set-TextValue 'E21' '  +2.64%  '

# Row 22
Set-TextValue 'D22' '388.39'
Set-TextValue 'E22' '  +5.19%  '

# Row 23
Set-TextValue 'E23' '  +0.14%  '

# Row 25
Set-TextValue 'D25' '70.72'
Set-TextValue 'E25' '  +1.97%  '

# Row 26
Set-TextValue 'E26' '  +10.64%  '

# Row 27
Set-TextValue 'D27' '0.180'
Set-TextValue 'E27' '  +6.20%  '

# Row 28
Set-TextValue 'E28' '  +18.22%  '

# Row 29
Set-TextValue 'E29' '  +0.16%  '

# Row 30
Set-TextValue 'E30' '  +7.67%  '

# Row 31
Set-TextValue 'D31' '6.47'
Set-TextValue 'E31' '  +5.91%  '

# Row 32
Set-TextValue 'D32' '23.07'
Set-TextValue 'E32' '  +2.70%  '

# Row 33
Set-TextValue 'E33' '  +10.23%  '

# Row 34
Set-TextValue 'E34' '  +5.42%  '

# Row 35
Set-TextValue 'D35' '6.73'
Set-TextValue 'E35' '  +2.47%  '

# Row 36
Set-TextValue 'E36' '  +9.28%  '

# Row 37
Set-TextValue 'D37' '158.58'
Set-TextValue 'E37' '  +1.35%  '

# Row 38
Set-TextValue 'B38' 'Stacks'
Set-TextValue 'C38' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '1.88'
Set-TextValue 'E38' '  +12.02%  '

# Row 39
Set-TextValue 'B39' 'EnergySwap'
Set-TextValue 'C39' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D39' '27.62'
Set-TextValue 'E39' '  +3.86%  '

# Row 40
Set-TextValue 'D40' '0.0750'
Set-TextValue 'E40' '  +6.21%  '

# Row 41
Set-TextValue 'D41' '2.877.14'
Set-TextValue 'E41' '  +2.89%  '

# Row 42
Set-TextValue 'E42' '  +4.22%  '

# Row 43
Set-TextValue 'E43' '  +5.80%  '

# Row 44
Set-TextValue 'D44' '41.15'

# Row 45
Set-TextValue 'D45' '4.29'
Set-TextValue 'E45' '  +0.87%  '

# Row 46
Set-TextValue 'E46' '  +7.83%  '

# Row 47
Set-TextValue 'D47' '3.434.49'
Set-TextValue 'E47' '  +6.23%  '

# Row 48
Set-TextValue 'D48' '22.07'
Set-TextValue 'E48' '  +6.84%  '

# Row 49
Set-TextValue 'D49' '298.62'
Set-TextValue 'E49' '  +12.76%  '

# Row 50
Set-TextValue 'B50' 'Stellar'
Set-TextValue 'C50' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D50' '0.103'
Set-TextValue 'E50' '  -0.84%  '

# Row 51
Set-TextValue 'B51' 'Cosmos'
Set-TextValue 'C51' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D51' '6.32'
Set-TextValue 'E51' '  +2.91%  '
